# "have advance use the same codes"
# The "Advance" rows (34-49) used to encode their VAR1 (column D) values with
# special one-off codes (AX, AOE, A). This edit makes them reuse the same
# codes already used by the "Play" rows above: X (new), E, FC, P.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("codes")

$ws.Range("D34").Value = "X"
$ws.Range("D37").Value = "E"
$ws.Range("D38").Value = "E"
$ws.Range("D39").Value = "E"
$ws.Range("D40").Value = "E"
$ws.Range("D41").Value = "FC"
$ws.Range("D42").Value = "P"
$ws.Range("D43").Value = "E"
$ws.Range("D46").Value = "E"
$ws.Range("D47").Value = "E"
$ws.Range("D48").Value = "FC"
$ws.Range("D49").Value = "FC"
